# Update cryptos list with latest price/volume data (GitHub Actions scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.032.61"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "3.120.92"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "593.83"
$ws.Range("D6").Value = "136.28"
$ws.Range("E6").Value = "  -5.53%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.118.57"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("E10").Value = "  -4.07%  "
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -5.46%  "
$ws.Range("D14").Value = "34.18"
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("D15").Value = "3.633.55"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "62.971.84"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").Value = "3.121.94"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").Value = "471.44"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "14.09"
$ws.Range("E21").Value = "  -4.05%  "
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "85.92"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").Value = "12.90"
$ws.Range("E25").Value = "  -4.28%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "7.88"
$ws.Range("E28").Value = "  -6.95%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "6.91"
$ws.Range("E29").Value = "  -5.00%  "
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "26.71"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  -6.17%  "
$ws.Range("E34").Value = "  -5.43%  "
$ws.Range("D35").Value = "1.08"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("D37").Value = "51.88"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").Value = "0.0₃0694"
$ws.Range("E38").Value = "  -11.37%  "
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("D40").Value = "419.27"
$ws.Range("E40").Value = "  -6.20%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "2.893.99"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  -12.99%  "
$ws.Range("E44").Value = "  -6.07%  "
$ws.Range("D45").Value = "0.264"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  -6.14%  "
$ws.Range("D48").Value = "25.47"
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("E50").Value = "  -7.69%  "
$ws.Range("D51").Value = "119.69"
$ws.Range("E51").Value = "  -0.36%  "
